$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.948.42"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.893.78"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").Value = "243.92"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.3136"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "25.75"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "0.07352"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "0.08073"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "0.7721"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "5.508"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").Value = "1.880.92"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "94.32"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "6.224"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "29.906.10"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "14.01"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "246.64"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "0.000007857"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "8.139"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.142.28"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "0.1582"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "162.35"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "2.031"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "1.424"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "4.478"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "4.067"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "0.7545"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").Value = "2.792"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "74.62"
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("D42").Value = "0.4480"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").Value = "1.101.80"
$ws.Range("E43").Value = "  +6.75%  "
$ws.Range("D44").Value = "6.004"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").Value = "0.8510"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D47").Value = "1.892"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("D48").Value = "102.41"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").Value = "7.549"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "9.787"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "2.999"
$ws.Range("E51").Value = "  +3.06%  "
